$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A39").Value = "Mattia Spagnolli "
$ws.Range("B39").Value = "Lorenzo Canali | CGB Gamberoni"
$ws.Range("C39").Value = "Filippo Benetti | I Magnifici"
$ws.Range("D39").Value = "Riccardo Zaffoni | U.SGUARNA"
$ws.Range("E39").Value = "Sebastiano Zoller | CGB Gamberoni"
$ws.Range("F39").Value = "Andrea Giordani | Clitoriders"
